# Separate suites for qa and stging and message update
#
# This updates the CreateAccount sheet's generated Selenium test e-mail
# addresses (E2:E6) and the stored automation password (F7), as well as
# the matching generated account ids on the Input sheet (U2:U6), so the
# workbook reflects a fresh batch of qa/staging test data.

$wb = $excel.ActiveWorkbook

# --- CreateAccount sheet -------------------------------------------------
$createAccount = $wb.Worksheets.Item("CreateAccount")

$createAccount.Range("E2").Value = "SeleniumKkkh@mailinator.com"
$createAccount.Range("E3").Value = "SeleniumZZrs@mailinator.com"
$createAccount.Range("E4").Value = "SeleniumgFkY@mailinator.com"
$createAccount.Range("E5").Value = "SeleniumgEdD@mailinator.com"
$createAccount.Range("E6").Value = "SeleniumSoRU@mailinator.com"

# Updated automation password used for the qa/staging suites
$createAccount.Range("F7").Value = "Automation1495!"

# Keep column E sized to fit the (slightly different) new e-mail text
$createAccount.Columns.Item(5).EntireColumn.AutoFit()

# --- Input sheet -----------------------------------------------------------
$inputSheet = $wb.Worksheets.Item("Input")

$inputSheet.Range("U2").Value = "51502122"
$inputSheet.Range("U3").Value = "51502129"
$inputSheet.Range("U4").Value = "51502143"
$inputSheet.Range("U5").Value = "51502146"
$inputSheet.Range("U6").Value = "51503402"
